$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 62
$ws.Range("H62").Value = 2500
$ws.Range("I62").Value = 2500
$ws.Range("K62").Value = 2500
$ws.Range("M62").Value = -1876
# row 65
$ws.Range("H65").Value = 2500
$ws.Range("I65").Value = 2500
$ws.Range("K65").Value = 12500
$ws.Range("M65").Value = -9380
# row 70
$ws.Range("H70").Value = 2020.6666
$ws.Range("I70").Value = 2134.4
$ws.Range("J70").Value = 1793.2
$ws.Range("K70").Value = 6403.200000000001
$ws.Range("L70").Value = 5379.6
$ws.Range("M70").Value = -6133.200000000001
$ws.Range("N70").Value = -5919.6
# row 73
$ws.Range("H73").Value = 2020.6666
$ws.Range("I73").Value = 2134.4
$ws.Range("J73").Value = 1793.2
$ws.Range("K73").Value = 6403.200000000001
$ws.Range("L73").Value = 5379.6
$ws.Range("M73").Value = -5467.200000000001
$ws.Range("N73").Value = -7251.6
# row 98
$ws.Range("H98").Value = 1519.8572
$ws.Range("I98").Value = 1519.8572
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 1519.8572
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -21.85719999999992
$ws.Range("N98").ClearContents()
# row 109
$ws.Range("H109").Value = 44930.668
$ws.Range("J109").Value = 44930.668
$ws.Range("L109").Value = 44930.668
$ws.Range("N109").Value = -47704.668
# row 112
$ws.Range("H112").Value = 2739.7273
$ws.Range("I112").Value = 1546.6666
$ws.Range("J112").Value = 3187.125
$ws.Range("K112").Value = 4639.9998
$ws.Range("L112").Value = 9561.375
$ws.Range("M112").Value = -3531.9998
$ws.Range("N112").Value = -11777.375
# row 122
$ws.Range("H122").Value = 1519.8572
$ws.Range("I122").Value = 1519.8572
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4559.571599999999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2109.571599999999
$ws.Range("N122").ClearContents()
# row 132
$ws.Range("H132").Value = 951.29034
$ws.Range("I132").Value = 1047.6923
$ws.Range("J132").Value = 450
$ws.Range("K132").Value = 3143.0769
$ws.Range("L132").Value = 1350
$ws.Range("M132").Value = -613.0769
$ws.Range("N132").Value = -6410
# row 135
$ws.Range("H135").Value = 1051.5
$ws.Range("I135").Value = 1051.5
$ws.Range("K135").Value = 9463.5
$ws.Range("M135").Value = -6928.5
# row 138
$ws.Range("H138").Value = 15749.75
$ws.Range("J138").Value = 15749.75
$ws.Range("L138").Value = 47249.25
$ws.Range("N138").Value = -57529.25

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 74
$ws.Range("H74").Value = 3331.6667
$ws.Range("I74").Value = 3331.6667
$ws.Range("K74").Value = 3331.6667
$ws.Range("M74").Value = -2457.6667
# row 77
$ws.Range("H77").Value = 3331.6667
$ws.Range("I77").Value = 3331.6667
$ws.Range("K77").Value = 16658.3335
$ws.Range("M77").Value = -12290.3335

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 81
$ws.Range("H81").Value = 33780
$ws.Range("J81").Value = 33780
$ws.Range("L81").Value = 33780
$ws.Range("N81").Value = -35902
# row 84
$ws.Range("H84").Value = 33780
$ws.Range("J84").Value = 33780
$ws.Range("L84").Value = 101340
$ws.Range("N84").Value = -111948
# row 99
$ws.Range("H99").Value = 1996.6666
$ws.Range("I99").Value = 1995
$ws.Range("K99").Value = 1995
$ws.Range("M99").Value = -497
# row 140
$ws.Range("H140").Value = 170000
$ws.Range("J140").Value = 170000
$ws.Range("L140").Value = 170000
$ws.Range("N140").Value = -180360

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 7
$ws.Range("H7").Value = 2621.5
$ws.Range("I7").Value = 2578.5
$ws.Range("K7").Value = 2578.5
$ws.Range("M7").Value = -2465.5
# row 99
$ws.Range("H99").Value = 6768.375
$ws.Range("I99").Value = 4500.7144
$ws.Range("J99").Value = 8532.111000000001
$ws.Range("K99").Value = 4500.7144
$ws.Range("L99").Value = 8532.111000000001
$ws.Range("M99").Value = -3002.7144
$ws.Range("N99").Value = -11528.111
# row 126
$ws.Range("H126").Value = 6768.375
$ws.Range("I126").Value = 4500.7144
$ws.Range("J126").Value = 8532.111000000001
$ws.Range("K126").Value = 13502.1432
$ws.Range("L126").Value = 25596.333
$ws.Range("M126").Value = -11032.1432
$ws.Range("N126").Value = -30536.333

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 5
$ws.Range("H5").Value = 703.2143
$ws.Range("I5").Value = 483.55554
$ws.Range("J5").Value = 1098.6
$ws.Range("K5").Value = 1450.66662
$ws.Range("L5").Value = 3295.8
$ws.Range("M5").Value = -1338.66662
$ws.Range("N5").Value = -3519.8
# row 135
$ws.Range("H135").Value = 703.2143
$ws.Range("I135").Value = 483.55554
$ws.Range("J135").Value = 1098.6
$ws.Range("K135").Value = 4351.99986
$ws.Range("L135").Value = 9887.4
$ws.Range("M135").Value = -1816.99986
$ws.Range("N135").Value = -14957.4

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 102
$ws.Range("H102").Value = 1136.3334
$ws.Range("I102").Value = 604
$ws.Range("K102").Value = 604
$ws.Range("M102").Value = 1018
# row 113
$ws.Range("H113").Value = 1250.7142
$ws.Range("I113").Value = 1250.7142
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1250.7142
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 919.2858000000001
$ws.Range("N113").ClearContents()
# row 126
$ws.Range("H126").Value = 74899464
$ws.Range("I126").Value = 124830040
$ws.Range("J126").Value = 3591.1667
$ws.Range("K126").Value = 374490120
$ws.Range("L126").Value = 10773.5001
$ws.Range("M126").Value = -374487650
$ws.Range("N126").Value = -15713.5001
# row 132
$ws.Range("H132").Value = 3972.4
$ws.Range("I132").Value = 3517.3333
$ws.Range("K132").Value = 10551.9999
$ws.Range("M132").Value = -8021.999899999999

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 40
$ws.Range("H40").Value = 10437.272
$ws.Range("I40").Value = 1495
$ws.Range("K40").Value = 1495
$ws.Range("M40").Value = -1359
# row 122
$ws.Range("H122").Value = 8618.375
$ws.Range("I122").Value = 8618.375
$ws.Range("K122").Value = 25855.125
$ws.Range("M122").Value = -23405.125
# row 132
$ws.Range("H132").Value = 9610.223
$ws.Range("I132").Value = 6916.1665
$ws.Range("K132").Value = 20748.4995
$ws.Range("M132").Value = -18218.4995

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 132
$ws.Range("H132").Value = 2987.1365
$ws.Range("I132").Value = 2214.625
$ws.Range("J132").Value = 5047.1665
$ws.Range("K132").Value = 6643.875
$ws.Range("L132").Value = 15141.4995
$ws.Range("M132").Value = -4113.875
$ws.Range("N132").Value = -20201.4995
